# Support for plans without sticky-ids or levels
# Rename the header row on the PV-Test-03 sheet to the new column names,
# and move the active selection from F4 to E2.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PV-Test-03")

$ws.Range("A1").Value = "Row ID"
$ws.Range("C1").Value = "Task"
$ws.Range("E1").Value = "Start Date"
$ws.Range("F1").Value = "End Date"

$ws.Range("E2").Select() | Out-Null
